$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "End Date" (F3) should have the same data type as "Start Date" (E3):
# change it from a date serial number to the text value used in E3.
# Copy E3's cell format to F3 first, then copy E3's value to F3, so
# that F3 ends up as a text cell (matching E3) instead of Excel
# re-interpreting the text as a date.
$ws.Range("E3").Copy()
$ws.Range("F3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E3").Copy()
$ws.Range("F3").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0

# Update "Total Absence Requested" (G3) to reflect the corrected dates.
$ws.Range("G3").Value = 1
